# Heatmap section now saves the epitables used to generate the heatmaps
#
# 1. "Average TP2 longitude" row: merge the two runs "tp2_avg_l" + "ong"
#    into a single run "tp2_avg_long".
# 2-5. Remove the (now-unwanted) spell-check proofErr bracketing around
#    four single-word column-name cells (actual_size_change, num_novs,
#    actual_growth_rate, new_growth).
# 6. "Type (1,2,3,4)" row: fill the previously-empty "Base" cell with "type".

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Merge the split run into one run -----------------------------
$d.Content.Find.Execute("tp2_avg_long", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "tp2_avg_long", 2)

# --- helper: rewrite a cell's single paragraph to drop proofErr wrap --
function Clear-ProofErr($cell, $text) {
    $para = $cell.Range.Paragraphs.Item(1)
    $para.Range.Delete()
    $cell.Range.Text = $text
}

# --- 2-5. Strip proofErr spellStart/spellEnd around the plain column
#          name cells -------------------------------------------------
Clear-ProofErr $t.Cell(39, 2) "actual_size_change"
Clear-ProofErr $t.Cell(41, 2) "num_novs"
Clear-ProofErr $t.Cell(42, 2) "actual_growth_rate"
Clear-ProofErr $t.Cell(43, 2) "new_growth"

# --- 6. Fill in the empty "Base" cell for the "Type (1,2,3,4)" row ----
$t.Cell(44, 2).Range.Text = "type"
